# Update cryptos list: price (D) and volume(1h) (E) columns for changed rows
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "23.521.05"
$ws.Range("E2").Value = "  +0.28%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.642.09"
$ws.Range("E3").Value = "  +0.10%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9992"
$ws.Range("E4").Value = "  -0.27%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.9994"
$ws.Range("E5").Value = "  -0.27%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "304.05"
$ws.Range("E6").Value = "  -0.09%  "
$ws.Range("E7").Value = "  +0.24%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "51.92"
$ws.Range("E8").Value = "  -0.74%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3628"
$ws.Range("E9").Value = "  -0.50%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.08191"
$ws.Range("E10").Value = "  +1.01%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.238"
$ws.Range("E11").Value = "  -0.91%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.9985"
$ws.Range("E12").Value = "  -0.33%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "22.61"
$ws.Range("E13").Value = "  -1.38%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.476"
$ws.Range("E14").Value = "  -2.48%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.382"
$ws.Range("E15").Value = "  +1.33%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.00001241"
$ws.Range("E16").Value = "  -0.94%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.634.33"
$ws.Range("E17").Value = "  -0.24%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "95.27"
$ws.Range("E18").Value = "  +1.17%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06947"
$ws.Range("E19").Value = "  +0.10%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.58"
$ws.Range("E20").Value = "  -3.23%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.574"
$ws.Range("E21").Value = "  +0.41%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.9999"
$ws.Range("E22").Value = "  -0.29%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "12.53"
$ws.Range("E23").Value = "  -2.58%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "23.512.77"
$ws.Range("E24").Value = "  +0.23%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.525"
$ws.Range("E25").Value = "  +2.72%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.082"
$ws.Range("E26").Value = "  -5.43%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "21.20"
$ws.Range("E27").Value = "  +0.02%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "152.35"
$ws.Range("E28").Value = "  +2.00%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "133.41"
$ws.Range("E30").Value = "  -1.82%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.817.16"
$ws.Range("E31").Value = "  -0.28%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.095"
$ws.Range("E32").Value = "  +13.29%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.627"
$ws.Range("E33").Value = "  -3.95%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.148"
$ws.Range("E34").Value = "  -7.22%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "11.54"
$ws.Range("E35").Value = "  +4.88%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02768"
$ws.Range("E36").Value = "  -3.71%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.2508"
$ws.Range("E37").Value = "  -2.10%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.08779"
$ws.Range("E38").Value = "  -1.36%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "6.035"
$ws.Range("E39").Value = "  -3.83%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.07067"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.7077"
$ws.Range("E41").Value = "  -0.73%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.356"
$ws.Range("E42").Value = "  -1.49%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "12.34"
$ws.Range("E43").Value = "  -1.96%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "15.65"
$ws.Range("E44").Value = "  -4.97%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.6567"
$ws.Range("E45").Value = "  +0.04%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.9993"
$ws.Range("E46").Value = "  -0.03%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.289"
$ws.Range("E47").Value = "  -3.03%  "
$ws.Range("E48").Value = "  -0.80%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.07990"
$ws.Range("E49").Value = "  -0.12%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "129.14"
$ws.Range("E50").Value = "  +1.16%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.199"
$ws.Range("E51").Value = "  -2.14%  "
